# Daily attendance processing - 2025-11-12 08:55:00
#
# This script updates the "Recorded By" column (column G) on the active
# worksheet: for the specific attendance rows listed below, the first two
# comma-separated entries in the "Recorded By" value are swapped (e.g.
# "dnasr281@gmail.com, System" becomes "System, dnasr281@gmail.com").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows (in column G) whose "Recorded By" value needs its first two
# comma-separated entries swapped.
$rowsToFix = @(
    2,3,4,6,7,10,11,12,13,14,15,17,18,19,20,21,22,24,26,28,29,30,
    32,33,36,37,38,39,40,41,43,44,45,46,47,48,50,52,54,55,56,58,59,
    62,63,64,65,66,67,69,70,71,72,73,74,76,78,83,84,85,86,90,92,93,
    94,96,99,101,109,110,111,112,116,118,119,120,122,125,127,135,
    136,137,138,142,144,145,146,148,151,153
)

foreach ($r in $rowsToFix) {
    $cell = $ws.Cells.Item($r, 7)   # Column G = 7 ("Recorded By")
    $current = [string]$cell.Value2

    # Split on comma, trim whitespace around each part.
    $parts = $current.Split(",") | ForEach-Object { $_.Trim() }

    if ($parts.Count -ge 2) {
        $first = $parts[0]
        $second = $parts[1]
        $rest = @()
        if ($parts.Count -gt 2) {
            $rest = $parts[2..($parts.Count - 1)]
        }

        $newParts = @($second, $first) + $rest
        $cell.Value = [string]::Join(", ", $newParts)
    }
}
